$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 1.88
$ws.Range("H2").Value2 = 4.8
$ws.Range("I2").Value2 = 4.9
$ws.Range("N2").Value2 = 3.35
$ws.Range("P2").Value2 = 1.82
$ws.Range("U2").Value2 = 1.94
$ws.Range("W2").Value2 = 2.12
$ws.Range("X2").Value2 = 13
$ws.Range("AB2").Value2 = 7.8

# Row 3
$ws.Range("F3").Value2 = 1.77
$ws.Range("G3").Value2 = 1.94
$ws.Range("H3").Value2 = 5
$ws.Range("I3").Value2 = 6.4
$ws.Range("J3").Value2 = 3.35
$ws.Range("K3").Value2 = 3.8
$ws.Range("L3").Value2 = 1.39
$ws.Range("M3").Value2 = 1.1
$ws.Range("N3").Value2 = 2.96
$ws.Range("O3").Value2 = 1.42
$ws.Range("P3").Value2 = 1.66
$ws.Range("Q3").Value2 = 2.22
$ws.Range("R3").Value2 = 1.24
$ws.Range("S3").Value2 = 4.3
$ws.Range("T3").Value2 = 2.06
$ws.Range("U3").Value2 = 1.78
$ws.Range("V3").Value2 = 1.18
$ws.Range("W3").Value2 = 2.06
$ws.Range("X3").Value2 = 13
$ws.Range("Y3").Value2 = 18
$ws.Range("Z3").Value2 = 50
$ws.Range("AB3").Value2 = 8.199999999999999
$ws.Range("AC3").Value2 = 9.6
$ws.Range("AD3").Value2 = 27
$ws.Range("AF3").Value2 = 11.5
$ws.Range("AG3").Value2 = 12.5
$ws.Range("AH3").Value2 = 29
$ws.Range("AJ3").Value2 = 23
$ws.Range("AK3").Value2 = 26
$ws.Range("AL3").Value2 = 60
$ws.Range("AN3").Value2 = 19

# Row 4
$ws.Range("F4").Value2 = 2.56
$ws.Range("G4").Value2 = 3.3
$ws.Range("H4").Value2 = 2.56
$ws.Range("I4").Value2 = 3.3
$ws.Range("J4").Value2 = 2.8
$ws.Range("K4").Value2 = 3.6
$ws.Range("L4").Value2 = 1.39
$ws.Range("M4").Value2 = 1.07
$ws.Range("N4").Value2 = 3.05
$ws.Range("O4").Value2 = 1.38
$ws.Range("P4").Value2 = 1.71
$ws.Range("Q4").Value2 = 2.02
$ws.Range("R4").Value2 = 1.27
$ws.Range("S4").Value2 = 3.55
$ws.Range("T4").Value2 = 1.81
$ws.Range("U4").Value2 = 1.96
$ws.Range("V4").Value2 = 1.44
$ws.Range("W4").Value2 = 1.43
$ws.Range("X4").Value2 = 990
$ws.Range("Y4").Value2 = 990
$ws.Range("AB4").Value2 = 12
$ws.Range("AC4").Value2 = 9

# Row 5
$ws.Range("H5").Value2 = 1.83
$ws.Range("I5").Value2 = 1.84
$ws.Range("J5").Value2 = 3.65
$ws.Range("K5").Value2 = 3.7
$ws.Range("L5").Value2 = 1.47
$ws.Range("T5").Value2 = 2.06
$ws.Range("V5").Value2 = 2.18
$ws.Range("Z5").Value2 = 10
$ws.Range("AM5").Value2 = 150
$ws.Range("AO5").Value2 = 14.5

# Row 6
$ws.Range("F6").Value2 = 2.44
$ws.Range("G6").Value2 = 3.55
$ws.Range("H6").Value2 = 2.26
$ws.Range("I6").Value2 = 3.15
$ws.Range("J6").Value2 = 3
$ws.Range("K6").Value2 = 7
$ws.Range("L6").Value2 = 1.28
$ws.Range("N6").Value2 = 2.14
$ws.Range("P6").Value2 = 1.94
$ws.Range("Q6").Value2 = 1.66
$ws.Range("R6").Value2 = 1.32
$ws.Range("S6").Value2 = 2.66
$ws.Range("V6").Value2 = 1.47
$ws.Range("W6").Value2 = 1.4

# Row 7
$ws.Range("F7").Value2 = 6
$ws.Range("G7").Value2 = 6.2
$ws.Range("H7").Value2 = 1.79
$ws.Range("I7").Value2 = 1.8
$ws.Range("J7").Value2 = 3.6
$ws.Range("K7").Value2 = 3.65
$ws.Range("Q7").Value2 = 2.42
$ws.Range("V7").Value2 = 2.24
$ws.Range("W7").Value2 = 1.19
$ws.Range("AB7").Value2 = 15
$ws.Range("AF7").Value2 = 44
$ws.Range("AG7").Value2 = 24
$ws.Range("AH7").Value2 = 26
$ws.Range("AL7").Value2 = 120

# Row 8
$ws.Range("H8").Value2 = 1.75
$ws.Range("I8").Value2 = 1.76
$ws.Range("J8").Value2 = 4
$ws.Range("K8").Value2 = 4.1
$ws.Range("L8").Value2 = 1.42
$ws.Range("N8").Value2 = 3.65
$ws.Range("O8").Value2 = 1.35
$ws.Range("P8").Value2 = 1.91
$ws.Range("S8").Value2 = 3.7
$ws.Range("T8").Value2 = 1.98
$ws.Range("V8").Value2 = 2.3
$ws.Range("AA8").Value2 = 17.5
$ws.Range("AE8").Value2 = 18.5
$ws.Range("AL8").Value2 = 80
$ws.Range("AO8").Value2 = 12

# Row 9
$ws.Range("F9").Value2 = 2.9
$ws.Range("G9").Value2 = 2.98
$ws.Range("H9").Value2 = 2.62
$ws.Range("I9").Value2 = 2.66
$ws.Range("L9").Value2 = 1.43
$ws.Range("N9").Value2 = 3.6
$ws.Range("O9").Value2 = 1.35
$ws.Range("Q9").Value2 = 2.06
$ws.Range("S9").Value2 = 3.8
$ws.Range("V9").Value2 = 1.6
$ws.Range("W9").Value2 = 1.51
$ws.Range("Y9").Value2 = 11
$ws.Range("Z9").Value2 = 16.5
$ws.Range("AA9").Value2 = 40
$ws.Range("AB9").Value2 = 12
$ws.Range("AF9").Value2 = 19.5
$ws.Range("AJ9").Value2 = 48
$ws.Range("AL9").Value2 = 55
$ws.Range("AM9").Value2 = 120

# Row 10
$ws.Range("F10").Value2 = 3.35
$ws.Range("G10").Value2 = 3.4
$ws.Range("H10").Value2 = 2.36
$ws.Range("I10").Value2 = 2.38
$ws.Range("J10").Value2 = 3.55
$ws.Range("K10").Value2 = 3.6
$ws.Range("T10").Value2 = 1.8
$ws.Range("V10").Value2 = 1.72
$ws.Range("Y10").Value2 = 10
$ws.Range("AE10").Value2 = 24
$ws.Range("AK10").Value2 = 34
$ws.Range("AL10").Value2 = 48

# Row 11
$ws.Range("F11").Value2 = 1.6
$ws.Range("G11").Value2 = 1.62
$ws.Range("J11").Value2 = 4.2
$ws.Range("P11").Value2 = 2.26
$ws.Range("S11").Value2 = 2.74
$ws.Range("U11").Value2 = 2.14
$ws.Range("W11").Value2 = 2.58
$ws.Range("X11").Value2 = 19.5
$ws.Range("Z11").Value2 = 75
$ws.Range("AE11").Value2 = 85
$ws.Range("AG11").Value2 = 10.5
$ws.Range("AH11").Value2 = 19
$ws.Range("AI11").Value2 = 85
$ws.Range("AM11").Value2 = 90
$ws.Range("AN11").Value2 = 7.4
$ws.Range("AO11").Value2 = 110
